# "added group 13 to sheet 11"
# Populate sheet "11" (sheet1 / rId1) with the four new names that were
# appended to the shared-string table, and make that sheet the active one
# (tabSelected / active tab / selection moves from sheet "12" to sheet "11").

$wb = $excel.ActiveWorkbook

$ws11 = $wb.Worksheets.Item("11")

$ws11.Range("A1").Value = "Bar Eckstien"
$ws11.Range("A2").Value = "Itai Orr"
$ws11.Range("A3").Value = "Ofri Efrati"
$ws11.Range("A4").Value = "Liav Teplizkiy"

# Make "11" the active sheet/tab and leave the selection on B14, matching
# the post-edit cursor position recorded in the workbook.
$ws11.Activate()
$ws11.Range("B14").Select()
